$d = $word.ActiveDocument

# 1) Update activation date
$d.Content.Find.Execute("Ativação: 01/01/2020", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2022", 2) | Out-Null

# 2) Trim trailing item 4 from Portuguese short summary ("Programa resumido")
$d.Content.Find.Execute(" 4) Trocadores de calor tubulares.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 3) Trim trailing item 4 from English short summary ("Programa resumido")
$d.Content.Find.Execute(" 4) Tubular heat exchangers.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 4) Trim trailing item 4 from Portuguese detailed program ("Programa")
$d.Content.Find.Execute(" 4) Determinação de coeficientes globais de troca de calor, balanços materiais e energéticos em trocadores tubulares do tipo casco e tubos.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 5) Fix double space before "coefficient" and trim trailing item 4 from English detailed program ("Programa")
$d.Content.Find.Execute("Diffusion  coefficient", $true, $false, $false, $false, $false, $true, 1, $false, "Diffusion coefficient", 2) | Out-Null
$d.Content.Find.Execute(" 4) Determination of overall heat transfer coefficients, material and energetic balances in shell-and-tube heat exchangers.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
